$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "updated date ranges for raw data"
# Column J holds "Mating Disrupt Start" dates; column K ("Mating Disrupt
# End" = J+90) recalculates automatically from the formula already in
# the sheet.
$ws.Range("J8").Value = 43634
$ws.Range("J9").Value = 43831

# Match the saved view state: the sheet was scrolled so row 4 is the
# top-most visible row, and the active selection moved from G5 to J9.
$ws.Range("J9").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
